$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9394612908363342
$ws.Range("B1").Value = 1.439923048019409
$ws.Range("C1").Value = 2.572733879089355
$ws.Range("D1").Value = 2.643383264541626
$ws.Range("E1").Value = 1.114669442176819
